# Update "想去人数" (interest count) figures across all four sheets to the
# values captured in the newer site snapshot (gh-pages output @ 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 299
$ws.Range("F4").Value = 1234
$ws.Range("F7").Value = 3817
$ws.Range("F9").Value = 743
$ws.Range("F10").Value = 1712
$ws.Range("F11").Value = 325
$ws.Range("F12").Value = 213
$ws.Range("F14").Value = 146
$ws.Range("F16").Value = 2064
$ws.Range("F17").Value = 320
$ws.Range("F21").Value = 219
$ws.Range("F22").Value = 14

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F9").Value = 124
$ws.Range("F10").Value = 88
$ws.Range("F12").Value = 79

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 303

# Sheet 4: 全部类型 (aggregate of the above three sheets)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 303
$ws.Range("F12").Value = 299
$ws.Range("F13").Value = 1234
$ws.Range("F19").Value = 3817
$ws.Range("F20").Value = 124
$ws.Range("F22").Value = 88
$ws.Range("F24").Value = 79
$ws.Range("F25").Value = 743
$ws.Range("F26").Value = 1712
$ws.Range("F27").Value = 325
$ws.Range("F29").Value = 213
$ws.Range("F31").Value = 146
$ws.Range("F34").Value = 2064
$ws.Range("F35").Value = 320
$ws.Range("F41").Value = 219
$ws.Range("F42").Value = 14
